# Build at 2022-09-26 16:07:08 UTC
# The source workbook had two "orphan" rows (13 & 14) that only carried a
# professor name in columns B/C with no label in column A. Those two rows
# are removed (all rows below shift up), and the paragraph-sized values
# that used to sit next to "Objetivos:"/"Programa resumido:"/"Programa:"/
# "Método:"/"Critério:"/"Norma de recuperação:"/"Bibliografia:" are swapped
# for the values that now line up after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two orphan rows (old rows 13 and 14); everything below shifts up.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()

# Row 10 ("Objetivos:") now shows the professor's name instead of the long
# course-objectives paragraph.
$ws.Range("B10").Value = "6634418 - Antonio Clelio Ribeiro"
$ws.Range("C10").Value = "6634418 - Antonio Clelio Ribeiro"

# Row 13 ("Programa resumido:") now shows the activation date value.
$ws.Range("B13").Value = "01/01/2019"
$ws.Range("C13").Value = "01/01/2019"

# Row 15 ("Programa:") now shows the professor's name.
$ws.Range("B15").Value = "6634418 - Antonio Clelio Ribeiro"
$ws.Range("C15").Value = "6634418 - Antonio Clelio Ribeiro"

# Row 18 ("Método:") now shows the second professor's name.
$ws.Range("B18").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C18").Value = "1285870 - Marcos Villela Barcza"

# Row 19 ("Critério:") now shows the teaching-method description.
$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula. discussão de castos práticos, visitas técnicas"
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula. discussão de castos práticos, visitas técnicas"

# Row 20 ("Norma de recuperação:") now shows the evaluation criteria text.
$ws.Range("B20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."

# Row 21 ("Bibliografia:") now shows the recovery-exam norm text.
$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
